# Update the encrypted "Посыл / Заповедь" strings in column B (rows 2-7)
# with freshly re-encrypted values (same plaintext, new ciphertext), and
# move the active selection from B10 to B11.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newB = "U2FsdGVkX187Z7DZ8yFC16RQzkiQV+stP4rB8oVF0Y2TWt5bNRPjiZk+HVhJ8AsLXUqKhbqoRBuWaSCj5QfG/telsAMBhdhhUI3h66j/it0g4O2elFlmo8o1vJvMDmQPEPuryqOLG7ly6knCdC76QlZrW62lh+3hwUoB0X8d0ku8kohl5nAd4Uy9k0qkvj+WkFzi3lEV9UsYQLfkiRF/ASFH6AFCrlxEwDDczx/vdNZg1vZxZJQzhzYwleSWQ1z2oCwpyvFAHp2sP4o8kJar434j4xCF3LFfWPGbQWZU8fansC/d+BfwsZVVvxZlFhs/RLd71yYalLGh1/auNVTyLjVZwJuxI68AKiuff33Uu9SIkEZdeeVt6PIa288yke80JN/xJeru4kfmxqXPnUpNUPj8QFVM1duTt/RnEg+SS6U237i46oCU+2KCnEDucruYPk/KbSd9AAmCT/GeeyqAuXyel9eSOsm1RKq7PW4Wn1LZQE9qYPG9gpeHKlmWwgDVPHqthl8Qmzv9YuZLpVeheEiEmGeww0U9x/3q7oyRrvgWVCcpH4whNm86E4CyW1efTRH3RguABuirEq1XrE2XdlM8dNCNF0w+zOOSjaTOL0sg7dWHqY5O5CTICifCuGdnqzf2RCT4jOScVG7K0d1nsBYADw1e3A8HNJ6cVRjXBA/FtOxlmdn421F/fz8nH6dmbqCZNDIMOHUPppAVdisUTANDuVuEzseXgcj6FSl40j5/3rPQx9G1cl+shgeLxHbcMPVdM93+I0CGREJIEMhnafjQGo+Dp3Gqs+hGH7mSUIgLe4RcdY4PX5lbcaAsAwIqpLtTG1+hkBUmVf1pokvlu76YRZA3Dg1SsXAY3EtLbPAnHeaLdFqQscmTFuFTZeJguraC9C3C3WbZ9ntHzec3EsRXRE2D44mP4SaPNnKICHvI7TwS2NJYgRv5ceavyQ6OMODOEr+2C9L82zJJ3qagV59N0LzIwGcLIdTEUA3zeGTX/sXySL/ONMFe0WDKwo8wtLabtJiCFcM0VZQfKzDrXwFpc9nYXNe9ShVm2u3YnfAInJBugHiRz3U1V8oNeKge1c5PEPUpu1cGkIxazKBsgH8KcYtv1RB3dJ020jIKLHFYO1hoWCEGqqaoCNX0/NKVTOhcqM4tmXrdl2i+P/N9yXr2UGwQU+ZEEy3gn86yNU88sEZCXeBDPhYop8+7BQPaC9Wv/DGbnNXq6Dgehdgd9KAM/Vr4aPrH3zU3p/UkSnO2seMu/u52TRWUv7a5JkyBVLy6uU/h9YZrOgTkTlr+K2ECy0gXLTGNp5OVPTHwPOk="
$newC = "U2FsdGVkX1+CGjVhD+CR0a49FzlkMhhI3K+Y8+C1QeyBumB2n3WxPoEqjJXAMwxSnSae2J24lQnenQKtpArmjqxj+t0uHXtQXVNaOqP1aQZ3HpZzllebjXdijZJcGrwiiQ+NiTw6zk7+WMoalu2Qfa8RyLWN65FkEmlAYQdVOka+cKSUoCmM6Gg+FVYCGedC5FPeKmVxtWAwLl9wto1zXTlbVUClpkQJ3u0FOoNEqvqFgG2r2YWbj+F7wpesQ5blKeAchqas2sz6FxPaDrcnA5alS9OR3YGp2UML/oq864LB8mDOykE8WsWlygLJ5ulC"

$ws.Range("B2").Value = $newB
$ws.Range("B3").Value = $newC
$ws.Range("B4").Value = $newB
$ws.Range("B5").Value = $newC
$ws.Range("B6").Value = $newB
$ws.Range("B7").Value = $newC

$ws.Range("B11").Select()
